# Update the win/transition-probability matrix on Sheet1 after simulating
# more games (see commit message: "added more games, sped up simulate game
# logic, and drafted optimization logic"). The underlying per-state game
# counts changed, which shifts several row probabilities; apply the new
# values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "P2"  = 0.2142857142857143
    "S2"  = 0.07142857142857142

    "P4"  = 0.5
    "S4"  = 0.5

    "D6"  = 0.03846153846153846
    "O6"  = 0.03846153846153846
    "Q6"  = 0.1153846153846154
    "S6"  = 0.4230769230769231

    "B7"  = 0.05263157894736842
    "F7"  = 0.05263157894736842
    "J7"  = 0.1052631578947368
    "O7"  = 0.05263157894736842
    "Q7"  = 0.1578947368421053
    "S7"  = 0.5789473684210527

    "B8"  = 0.03125
    "F8"  = 0.0625
    "J8"  = 0.0625
    "O8"  = 0.0625
    "Q8"  = 0.15625
    "R8"  = 0.09375
    "S8"  = 0.53125

    "B9"  = 0.1
    "F9"  = 0.1
    "J9"  = 0.1
    "Q9"  = 0.1
    "S9"  = 0.6

    "B10" = 0.1153846153846154
    "D10" = 0.01282051282051282
    "F10" = 0.1153846153846154
    "J10" = 0.08974358974358974
    "O10" = 0.01282051282051282
    "Q10" = 0.217948717948718
    "R10" = 0.03846153846153846
    "S10" = 0.3974358974358974

    "F11" = 0.03703703703703703
    "G11" = 0.1111111111111111
    "J11" = 0.07407407407407407
    "K11" = 0.1481481481481481
    "L11" = 0.5925925925925926
    "S11" = 0.03703703703703703

    "G12" = 0.7058823529411765
    "J12" = 0.1176470588235294
    "L12" = 0.05882352941176471
    "S12" = 0.1176470588235294

    "G13" = 0.5714285714285714
    "J13" = 0.2857142857142857
    "S13" = 0.1428571428571428

    "F15" = 0.08333333333333333
    "H15" = 0.08333333333333333
    "J15" = 0.25
    "K15" = 0.08333333333333333
    "O15" = 0.08333333333333333
    "S15" = 0.4166666666666667

    "H16" = 0.1666666666666667
    "I16" = 0.08333333333333333
    "J16" = 0.1666666666666667
    "K16" = 0.4166666666666667
    "M16" = 0.08333333333333333
    "S16" = 0.08333333333333333

    "F17" = 0.0625
    "H17" = 0.15625
    "I17" = 0.125
    "J17" = 0.21875
    "K17" = 0.03125
    "M17" = 0.0625
    "O17" = 0.0625
    "S17" = 0.28125

    "H18" = 0.1666666666666667
    "M18" = 0.1666666666666667
    "S18" = 0.1666666666666667

    "F19" = 0.03738317757009346
    "H19" = 0.205607476635514
    "I19" = 0.04672897196261682
    "J19" = 0.3644859813084112
    "K19" = 0.1588785046728972
    "M19" = 0.02803738317757009
    "O19" = 0.02803738317757009
    "S19" = 0.1308411214953271
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
